$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix "Objetivos:" content (row 10) ---
$objectives = "Proporcionar ao aluno conhecimento básico e compreensão de cinemática e dinâmica do corpo rígido. Desenvolver algumas aplicações práticas com ênfase em problemas bidimensionais. Apresentar conceitos fundamentais e exemplos das vibrações mecânicas."
$ws.Range("B10").Value = $objectives
$ws.Range("C10").Value = $objectives

# --- Insert a new row 13 to hold the professor name under "Docentes responsáveis:" ---
$ws.Rows.Item(13).Insert()
$ws.Range("A13").Clear()

# Copy number/alignment/font formatting from the existing B/C "value" columns (row 9)
# onto the freshly inserted row 13 cells, so they end up with the same style indices
# (s="2" for column B, s="3" for column C) as every other content row.
$ws.Range("B9").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C9").Copy()
$ws.Range("C13").PasteSpecial(-4122)

$prof = "7797767 - Viktor Pastoukhov"
$ws.Range("B13").Value = $prof
$ws.Range("C13").Value = $prof

# --- Fix "Programa resumido:" content (now row 14 after the insert above) ---
$shortSyllabus = "Cinemática do corpo rígidoDinâmica do pontoDinâmica do corpo rígido Introdução às vibrações mecânicas"
$ws.Range("B14").Value = $shortSyllabus
$ws.Range("C14").Value = $shortSyllabus

# --- Fix "Programa:" content (now row 16 after the insert above) ---
$syllabus = "Cinemática do corpo rígido:Aceleração e velocidade angulares. Vínculo e cinemática do corpo rígido. Rotação em torno de um eixo fixo. Movimento plano e centro de rotação. Composição de movimentos. Composição de movimentos de rotação.Dinâmica do ponto:Princípios da dinâmica do ponto. Teorema da resultante. Teorema da energia cinética para partícula. Teorema da quantidade de movimento.Dinâmica do corpo rígido:Teorema do movimento do baricentro. Teorema da energia cinética para um sistema de partículas. Teorema do momento angular para um sistema de partículas. Teorema da energia cinética para o corpo rígido. Teorema do momento angular para corpo rígido Exercícios de aplicação: problemas bidimensionais. Rotação do corpo rígido, Balanceamento. Movimento de um giroscópio.Introdução às vibrações mecânicas:Vibrações de sistemas mecânicos com um grau de liberdade: livres sem amortecimento, livres com amortecimento, forçadas. Vibrações de sistemas mecânicos com dois e mais graus de liberdade. Exemplos."
$ws.Range("B16").Value = $syllabus
$ws.Range("C16").Value = $syllabus

# --- Fix "Método:" content (now row 19 after the insert above) ---
$method = "A avaliação será composta por duas provas (P1 e P2)."
$ws.Range("B19").Value = $method
$ws.Range("C19").Value = $method

# --- Fix "Critério:" content (now row 20 after the insert above) ---
$criteria = "NS = NP1+NP2; NP1: questões da P1 valendo até 4p. no total; NP2: questões da P2 valendo até 6 p. no total."
$ws.Range("B20").Value = $criteria
$ws.Range("C20").Value = $criteria

# --- Fix "Norma de recuperação:" content (now row 21 after the insert above) ---
$recovery = "A recuperação consistirá de uma prova de Recuperação (R), que irá compor a nota final (NF) da seguinte forma: NF = (R + NS)/2."
$ws.Range("B21").Value = $recovery
$ws.Range("C21").Value = $recovery

# --- Fix "Bibliografia:" content (now row 22 after the insert above) ---
$biblio = "HIBBELER, R.C. Dinâmica - Mecânica para Engenharia. São Paulo: Pearson Brasil, 2011, 12ª ed., 608p. ISBN: 8576058146.BEER, F.P., JOHNSTON Jr., E.R., CLAUSEN, W. E., Mecânica Vetorial para Engenheiros - Dinâmica, 7ª Edição, McGraw-Hill, São Paulo, 2006, 1355 p. FRANÇA, L. N. F., MATSUMURA, A. Z. Mecânica Geral. Edgard Blücher, 2001, 235 p.SOTELO JR., J., FRANÇA, L.N.F., Introdução às vibrações mecânicas, Edgard Blücher, 2006, 168 p. ISBN: 9788521203384.GREENWOOD, D. T. Principles of Dynamics. New York: Prentice-Hall, 2nd ed, 1988, 552 p.TENENBAUM, R. A. Dinâmica. Editora UFRJ, 1997, 756 p.GIACAGLIA, G. E., Mecânica Geral, Editora Campus, Rio de Janeiro, 1982."
$ws.Range("B22").Value = $biblio
$ws.Range("C22").Value = $biblio
